$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two hyperlinks that used to live in column D (poc_email), along
# with their underlying relationships.
$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("D3").Hyperlinks.Delete()

# --- Re-author the table -------------------------------------------------
# Row 1: new 3-column header (POC_name / POC_designation / POC_contact)
$ws.Range("A1").Value = "POC_name"
$ws.Range("B1").Value = "POC_designation"
$ws.Range("C1").Value = "POC_contact"

# Row 3: new POC - Aryan
$ws.Range("A3").Value = "Aryan"
$ws.Range("B3").Value = "Deputy Coordinator"
$ws.Range("C3").Value = 1000000000

# Row 4: new POC - Aditya
$ws.Range("A4").Value = "Aditya"
$ws.Range("B4").Value = "Head"
$ws.Range("C4").Value = 1000000000

# Row 2: Faisal's designation updated last
$ws.Range("B2").Value = "Placement Coordinator"

# Drop the now-unused poc_email / department columns (D:E) entirely -
# clear their values and formatting (the hyperlink style lived in D).
$ws.Range("D1:E4").Clear()

# --- Cosmetic / view changes ---------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 19.42578125

$ws.Range("B5").Select() | Out-Null
